$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.784.61"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.854.87"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.60%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.79"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4404"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3814"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07419"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8870"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.57"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.860.19"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.520"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.720"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07193"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "85.28"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.039"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009085"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.032"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.53"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.799.46"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.285"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.28"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.089.20"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.18%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +6.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.26"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.000"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.358"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.10"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09083"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.215"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7740"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.012"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.586"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.034"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.152"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.56%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05295"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5194"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1672"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.888"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.760"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "110.83"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.78"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.035"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06585"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.713"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4717"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.898"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.44%  "
